# Auto-generated Excel COM-interop script to update cryptos worksheet data
# Updates Price (D) and Volume(1h) (E) columns, and two pairs of swapped rows
# (NEARProtocol/ImmutableX at rows 33-34, Hedera/Mantle at rows 45-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.997.43"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.687.63"
$ws.Range("E3").Value = "  -1.73%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "3.163.96"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "62.912.51"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "2.690.67"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "346.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.512"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.05%  "
$ws.Range("D29").Value = "0.0₃0862"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "360.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.964"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.621"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0974"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.57%  "
